# "fix supplier reviews comments"
# Update the text of three requirement cells in column B of Sheet1:
#  - B27: supplier "add product" requirement now calls out IOS & Android platforms
#  - B30: Product Id naming convention drops the "[Color]" segment
#  - B31: mandatory/optional product-data wording reworded (photo moves to mandatory)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B27").Value = "Supplier can add product with data (product id,product photo, product price, product version, product platform [IOS &Android] )"
$ws.Range("B30").Value = "Product Id follow ID convention [Brand/Category]_[ProductName]_[Size/Version]"
$ws.Range("B31").Value = "If the supplier wants to add a new product the product data which is mandatory (product id, product price, product version, product photo& product platform)"

# Cursor/selection ends up on the last-reviewed cell, matching the saved view state.
[void]$ws.Range("B27").Select()
